# Roll the income-statement table forward by one reporting period:
#  - drop the oldest "12 ماهه منتهی به 1396/12" column content, shift each
#    period/date header left by one column, and introduce the new
#    1401/12 period (with its publish date) in column H
#  - shift every data row's values left by one column and compute the new
#    latest-period (column H) figure
#  - row 15 (هزینه کاهش ارزش دریافتنی‌ها) had a literal "-" placeholder in
#    column D; that now becomes a real 0 now that the period has rolled
#  - row 26 (سرمایه / capital) does not follow the simple shift pattern,
#    its raw values are rewritten explicitly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8: period headers ----
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# ---- Row 9: publish dates ----
$ws.Range("D9").Value = "1399-04-11 (7)"
$ws.Range("E9").Value = "1400-04-09 (8)"
$ws.Range("F9").Value = "1401-04-01 (8)"
$ws.Range("G9").Value = "1402-02-30 (8)"
$ws.Range("H9").Value = "1402-02-30 (2)"

# ---- Row 11: فروش (sales) ----
$ws.Range("D11").Value = 1544869
$ws.Range("E11").Value = 2185762
$ws.Range("F11").Value = 3567949
$ws.Range("G11").Value = 5958100
$ws.Range("H11").Value = 12813639

# ---- Row 12: بهای تمام شده کالای فروش رفته (COGS) ----
$ws.Range("D12").Value = -1296594
$ws.Range("E12").Value = -1742678
$ws.Range("F12").Value = -2787749
$ws.Range("G12").Value = -4663492
$ws.Range("H12").Value = -8536412

# ---- Row 13: سود (زیان) ناخالص (gross profit) ----
$ws.Range("D13").Value = 248275
$ws.Range("E13").Value = 443084
$ws.Range("F13").Value = 780200
$ws.Range("G13").Value = 1294608
$ws.Range("H13").Value = 4277227

# ---- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ----
$ws.Range("D14").Value = -72899
$ws.Range("E14").Value = -102267
$ws.Range("F14").Value = -177605
$ws.Range("G14").Value = -307324
$ws.Range("H14").Value = -892802

# ---- Row 15: هزینه کاهش ارزش دریافتنی‌ها (impairment expense) ----
# D15 used to hold the literal placeholder "-"; it now rolls to a real 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# ---- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (other operating income/exp) ----
$ws.Range("D16").Value = -7789
$ws.Range("E16").Value = -8646
$ws.Range("F16").Value = -27586
$ws.Range("G16").Value = -14591
$ws.Range("H16").Value = -8591

# ---- Row 17: سود (زیان) عملیاتی (operating profit) ----
$ws.Range("D17").Value = 167587
$ws.Range("E17").Value = 332171
$ws.Range("F17").Value = 575009
$ws.Range("G17").Value = 972693
$ws.Range("H17").Value = 3375834

# ---- Row 18: هزینه های مالی (finance costs) ----
$ws.Range("D18").Value = -9237
$ws.Range("E18").Value = -23437
$ws.Range("F18").Value = -5826
$ws.Range("G18").Value = -11698
$ws.Range("H18").Value = -108052

# ---- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (other non-op income/exp) ----
$ws.Range("D19").Value = 2155
$ws.Range("E19").Value = 40280
$ws.Range("F19").Value = 12025
$ws.Range("G19").Value = 50929
$ws.Range("H19").Value = 48612

# ---- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (pre-tax profit) ----
$ws.Range("D20").Value = 160505
$ws.Range("E20").Value = 349014
$ws.Range("F20").Value = 581208
$ws.Range("G20").Value = 1011924
$ws.Range("H20").Value = 3316394

# ---- Row 21: مالیات (tax) ----
$ws.Range("D21").Value = -37825
$ws.Range("E21").Value = -46703
$ws.Range("F21").Value = -101087
$ws.Range("G21").Value = -146985
$ws.Range("H21").Value = -398868

# ---- Row 22: سود (زیان) خالص عملیات در حال تداوم (net profit from continuing ops) ----
$ws.Range("D22").Value = 122680
$ws.Range("E22").Value = 302311
$ws.Range("F22").Value = 480121
$ws.Range("G22").Value = 864939
$ws.Range("H22").Value = 2917526

# ---- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (discontinued ops) - all zero, unchanged ----
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# ---- Row 24: سود (زیان) خالص (net profit) ----
$ws.Range("D24").Value = 122680
$ws.Range("E24").Value = 302311
$ws.Range("F24").Value = 480121
$ws.Range("G24").Value = 864939
$ws.Range("H24").Value = 2917526

# ---- Row 25: سود هر سهم پس از کسر مالیات (EPS after tax) ----
$ws.Range("D25").Value = 27
$ws.Range("E25").Value = 1008
$ws.Range("F25").Value = 107
$ws.Range("G25").Value = 193
$ws.Range("H25").Value = 651

# ---- Row 26: سرمایه (capital) - does not follow the simple left-shift ----
$ws.Range("D26").Value = 4484000
$ws.Range("E26").Value = 300000
$ws.Range("F26").Value = 4484000
$ws.Range("G26").Value = 4484000
$ws.Range("H26").Value = 4484000

# ---- Row 27: سود هر سهم بر اساس آخرین سرمایه (EPS based on latest capital) ----
$ws.Range("D27").Value = 27
$ws.Range("E27").Value = 67
$ws.Range("F27").Value = 107
$ws.Range("G27").Value = 193
$ws.Range("H27").Value = 651
